# Rotate the data among rows 2, 5, 7, 9:
#   new row2 <- old row5
#   new row5 <- old row7
#   new row7 <- old row9
#   new row9 <- old row2
#
# Columns that actually differ between the four rows (others are identical
# across all of them): A, B, E, F, G, H, K, L, M, N, Q, R, AC.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture "before" values for the 4 rows that are about to rotate ---
$rows = @(2, 5, 7, 9)

$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{
        A  = $ws.Range("A$r").Value2
        B  = $ws.Range("B$r").Value2
        E  = $ws.Range("E$r").Value2
        F  = $ws.Range("F$r").Value2
        G  = $ws.Range("G$r").Value2
        H  = $ws.Range("H$r").Value2
        K  = $ws.Range("K$r").Value2
        L  = $ws.Range("L$r").Value2
        M  = $ws.Range("M$r").Value2
        N  = $ws.Range("N$r").Value2
        Q  = $ws.Range("Q$r").Value2
        R  = $ws.Range("R$r").Value2
        AC = $ws.Range("AC$r").Value2
    }
}

# Source row for each destination row (cyclic rotation).
$srcFor = @{ 2 = 5; 5 = 7; 7 = 9; 9 = 2 }

foreach ($dst in $rows) {
    $src = $srcFor[$dst]
    $vals = $before[$src]

    $ws.Range("A$dst").Value2 = $vals.A
    $ws.Range("B$dst").Value2 = $vals.B
    $ws.Range("E$dst").Value2 = $vals.E
    $ws.Range("F$dst").Value2 = $vals.F
    $ws.Range("G$dst").Value2 = $vals.G
    $ws.Range("H$dst").Value2 = $vals.H

    if ([string]::IsNullOrEmpty($vals.K)) {
        $ws.Range("K$dst").ClearContents()
    } else {
        $ws.Range("K$dst").Value2 = $vals.K
    }
    if ([string]::IsNullOrEmpty($vals.L)) {
        $ws.Range("L$dst").ClearContents()
    } else {
        $ws.Range("L$dst").Value2 = $vals.L
    }
    if ([string]::IsNullOrEmpty($vals.M)) {
        $ws.Range("M$dst").ClearContents()
    } else {
        $ws.Range("M$dst").Value2 = $vals.M
    }
    if ([string]::IsNullOrEmpty($vals.N)) {
        $ws.Range("N$dst").ClearContents()
    } else {
        $ws.Range("N$dst").Value2 = $vals.N
    }

    $ws.Range("Q$dst").Value2 = $vals.Q
    $ws.Range("R$dst").Value2 = $vals.R

    if ([string]::IsNullOrEmpty($vals.AC)) {
        $ws.Range("AC$dst").ClearContents()
    } else {
        $ws.Range("AC$dst").Value2 = $vals.AC
    }
}
